# New crime data collected - weekly CompStat update
# Precinct volume/issue number and the reporting week dates advance by one week,
# and the weekly/28-day/YTD/2-year crime figures (and their derived % change
# columns) are refreshed with the newly collected counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: issue number and reporting week dates --------------------
# These cells hold "rich text" (several runs with identical formatting), so a
# plain text assignment collapses them to a single run - the rendered text is
# unchanged in appearance.
$ws.Range("A8").Value = "Volume 30   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# --- Cells that flip from a numeric "1"/"0" to the text placeholder "0" or
# "***.*" (meaning: no complaints / figure suppressed this period). A literal
# Value assignment of a digit-only string gets auto-coerced to a number by
# Excel, so we instead enter it as a text formula result, paste-special the
# *values* (collapses the formula to a plain shared-string literal, exactly
# like the other placeholder cells in the sheet) and finally paste-special
# the *formats* from a neighboring cell that already carries the same
# placeholder + style, so no new number-format style gets fabricated.
$ws.Range("C14").Formula = '="0"'
$ws.Range("C14").Copy()
$ws.Range("C14").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("D20").Formula = '="0"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Formula = '="***.*"'
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("C28").Formula = '="0"'
$ws.Range("C28").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("C29").Formula = '="0"'
$ws.Range("C29").Copy()
$ws.Range("C29").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Cell that flips from the text placeholder "0" back to a real number ---
$ws.Range("F30").Value = 1
$ws.Range("F30").NumberFormat = "#,##0"

# --- Refreshed weekly / 28-day / YTD / 2-year counts and % changes ---------
$ws.Range("F15").Value = 1
$ws.Range("N15").Value = -62.857142857142

$ws.Range("D16").Value = 15
$ws.Range("E16").Value = -40
$ws.Range("G16").Value = 68
$ws.Range("H16").Value = -52.941176470588
$ws.Range("I16").Value = 380
$ws.Range("J16").Value = 517
$ws.Range("K16").Value = -26.499032882011
$ws.Range("L16").Value = 12.759643916913
$ws.Range("M16").Value = 196.875
$ws.Range("N16").Value = -81.695568400770

$ws.Range("C17").Value = 6
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = -26.666666666666
$ws.Range("I17").Value = 404
$ws.Range("J17").Value = 383
$ws.Range("K17").Value = 5.483028720626
$ws.Range("L17").Value = 9.782608695652
$ws.Range("M17").Value = 164.052287581699
$ws.Range("N17").Value = -26.411657559198

$ws.Range("C18").Value = 12
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 42
$ws.Range("H18").Value = -35.714285714285
$ws.Range("I18").Value = 336
$ws.Range("J18").Value = 545
$ws.Range("K18").Value = -38.348623853211
$ws.Range("L18").Value = -4.815864022662
$ws.Range("M18").Value = 21.739130434782
$ws.Range("N18").Value = -84.615384615384

$ws.Range("C19").Value = 48
$ws.Range("D19").Value = 47
$ws.Range("E19").Value = 2.127659574468
$ws.Range("F19").Value = 158
$ws.Range("G19").Value = 186
$ws.Range("H19").Value = -15.053763440860
$ws.Range("I19").Value = 1856
$ws.Range("J19").Value = 1832
$ws.Range("K19").Value = 1.310043668122
$ws.Range("L19").Value = 70.745170193192
$ws.Range("M19").Value = 1.922020867655
$ws.Range("N19").Value = -75.967888126375

$ws.Range("I20").Value = 57
$ws.Range("K20").Value = 11.764705882352
$ws.Range("L20").Value = 29.545454545454
$ws.Range("M20").Value = 171.428571428571
$ws.Range("N20").Value = -80.479452054794

$ws.Range("C21").Value = 76
$ws.Range("D21").Value = 82
$ws.Range("E21").Value = -7.317073170731
$ws.Range("F21").Value = 256
$ws.Range("G21").Value = 344
$ws.Range("H21").Value = -25.581395348837
$ws.Range("I21").Value = 3049
$ws.Range("J21").Value = 3352
$ws.Range("K21").Value = -9.039379474940
$ws.Range("L21").Value = 38.213961922030
$ws.Range("M21").Value = 26.619601328903
$ws.Range("N21").Value = -76.305564190239

$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -33.333333333333
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 19
$ws.Range("H22").Value = -47.368421052631
$ws.Range("I22").Value = 161
$ws.Range("J22").Value = 159
$ws.Range("K22").Value = 1.257861635220
$ws.Range("L22").Value = 29.838709677419
$ws.Range("M22").Value = 38.793103448275

$ws.Range("C24").Value = 72
$ws.Range("D24").Value = 75
$ws.Range("E24").Value = -4
$ws.Range("F24").Value = 285
$ws.Range("G24").Value = 283
$ws.Range("H24").Value = 0.706713780918
$ws.Range("I24").Value = 3266
$ws.Range("J24").Value = 2694
$ws.Range("K24").Value = 21.232368225686
$ws.Range("L24").Value = 84.415584415584
$ws.Range("M24").Value = -14.903595622720

$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 56.25
$ws.Range("F25").Value = 81
$ws.Range("G25").Value = 66
$ws.Range("H25").Value = 22.727272727272
$ws.Range("I25").Value = 872
$ws.Range("J25").Value = 718
$ws.Range("K25").Value = 21.448467966573
$ws.Range("L25").Value = 22.128851540616
$ws.Range("M25").Value = 81.288981288981

$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0

$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = -28.571428571428
$ws.Range("F27").Value = 15
$ws.Range("G27").Value = 18
$ws.Range("H27").Value = -16.666666666666
$ws.Range("I27").Value = 179
$ws.Range("J27").Value = 185
$ws.Range("K27").Value = -3.243243243243
$ws.Range("L27").Value = 43.2

$ws.Range("H30").Value = -50
$ws.Range("I30").Value = 10
$ws.Range("K30").Value = -50
$ws.Range("L30").Value = -69.696969696969
